# Apply scheduled-runner profit/price updates to the Kraken_Profits workbook.
# Each leve row H:N block (currentAveragePrice/NQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# is refreshed from the latest market data snapshot. Blank cells in the source
# data (no NQ or HQ market) are cleared rather than written as 0.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5000
$ws.Range("J43").Value = 5000
$ws.Range("L43").Value = 5000
$ws.Range("N43").Value = -5138
$ws.Range("H58").Value = 1500
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = ""
$ws.Range("H70").Value = 3500
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""
$ws.Range("H73").Value = 3500
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""
$ws.Range("H80").Value = 865.44446
$ws.Range("I80").Value = 599.6667
$ws.Range("K80").Value = 1799.0001
$ws.Range("M80").Value = -801.0001
$ws.Range("H83").Value = 865.44446
$ws.Range("I83").Value = 599.6667
$ws.Range("K83").Value = 5397.0003
$ws.Range("M83").Value = -405.0002999999997
$ws.Range("H96").Value = 479.53845
$ws.Range("I96").Value = 369.9091
$ws.Range("J96").Value = 1082.5
$ws.Range("K96").Value = 1109.7273
$ws.Range("L96").Value = 3247.5
$ws.Range("M96").Value = 263.2727
$ws.Range("N96").Value = -5993.5
$ws.Range("H138").Value = 3603.95
$ws.Range("I138").Value = 3781
$ws.Range("J138").Value = 3426.9
$ws.Range("K138").Value = 11343
$ws.Range("L138").Value = 10280.7
$ws.Range("M138").Value = -6203
$ws.Range("N138").Value = -20560.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 100
$ws.Range("K12").Value = 100
$ws.Range("M12").Value = 68
$ws.Range("H86").Value = 2937.6
$ws.Range("I86").Value = 3666.6667
$ws.Range("J86").Value = 1844
$ws.Range("K86").Value = 3666.6667
$ws.Range("L86").Value = 1844
$ws.Range("M86").Value = -2543.6667
$ws.Range("N86").Value = -4090
$ws.Range("H89").Value = 2937.6
$ws.Range("I89").Value = 3666.6667
$ws.Range("J89").Value = 1844
$ws.Range("K89").Value = 18333.3335
$ws.Range("L89").Value = 9220
$ws.Range("M89").Value = -12717.3335
$ws.Range("N89").Value = -20452
$ws.Range("H94").Value = 2034.5
$ws.Range("I94").Value = 2038.3334
$ws.Range("K94").Value = 2038.3334
$ws.Range("M94").Value = -1587.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 344.13333
$ws.Range("I7").Value = 374
$ws.Range("J7").Value = 310
$ws.Range("K7").Value = 374
$ws.Range("L7").Value = 310
$ws.Range("M7").Value = -261
$ws.Range("N7").Value = -536
$ws.Range("H16").Value = 766.3333
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = ""
$ws.Range("H58").Value = 2397.5
$ws.Range("I58").Value = 2108.3333
$ws.Range("K58").Value = 2108.3333
$ws.Range("M58").Value = -1905.3333
$ws.Range("H86").Value = 4000
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -6246
$ws.Range("H89").Value = 4000
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -31232
$ws.Range("H113").Value = 766.3333
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""
$ws.Range("H122").Value = 2498.75
$ws.Range("I122").Value = 2666.1667
$ws.Range("J122").Value = 1996.5
$ws.Range("K122").Value = 7998.500100000001
$ws.Range("L122").Value = 5989.5
$ws.Range("M122").Value = -5548.500100000001
$ws.Range("N122").Value = -10889.5
$ws.Range("H134").Value = 912.5
$ws.Range("I134").Value = 912.5
$ws.Range("K134").Value = 2737.5
$ws.Range("M134").Value = -202.5
$ws.Range("H136").Value = 2397.5
$ws.Range("I136").Value = 2108.3333
$ws.Range("K136").Value = 6324.999899999999
$ws.Range("M136").Value = -3774.999899999999
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""
$ws.Range("H141").Value = 747078
$ws.Range("J141").Value = 1076219
$ws.Range("L141").Value = 1076219
$ws.Range("N141").Value = -1086579

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 8025
$ws.Range("I99").Value = 8025
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 24075
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -21829
$ws.Range("N99").Value = ""
$ws.Range("H109").Value = 896.4
$ws.Range("I109").Value = 870.5
$ws.Range("K109").Value = 2611.5
$ws.Range("M109").Value = -1571.5
$ws.Range("H113").Value = 897
$ws.Range("I113").Value = 1193.5
$ws.Range("J113").Value = 659.8
$ws.Range("K113").Value = 3580.5
$ws.Range("L113").Value = 1979.4
$ws.Range("M113").Value = -1410.5
$ws.Range("N113").Value = -6319.4
$ws.Range("H131").Value = 3276.0588
$ws.Range("I131").Value = 1899.2222
$ws.Range("J131").Value = 4825
$ws.Range("K131").Value = 5697.6666
$ws.Range("L131").Value = 14475
$ws.Range("M131").Value = -657.6665999999996
$ws.Range("N131").Value = -24555
$ws.Range("H138").Value = 4925
$ws.Range("I138").Value = 4566.6665
$ws.Range("K138").Value = 13699.9995
$ws.Range("M138").Value = -8559.999500000002
$ws.Range("H139").Value = 2337.5
$ws.Range("I139").Value = 1783.3334
$ws.Range("K139").Value = 5350.0002
$ws.Range("M139").Value = -210.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3590.3157
$ws.Range("I102").Value = 3512
$ws.Range("K102").Value = 3512
$ws.Range("M102").Value = -1890

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5155.857
$ws.Range("I61").Value = 7166.3335
$ws.Range("J61").Value = 3648
$ws.Range("K61").Value = 7166.3335
$ws.Range("L61").Value = 3648
$ws.Range("M61").Value = -6964.3335
$ws.Range("N61").Value = -4052
$ws.Range("H68").Value = 2565.1538
$ws.Range("I68").Value = 2589.7
$ws.Range("K68").Value = 2589.7
$ws.Range("M68").Value = -1840.7
$ws.Range("H71").Value = 2565.1538
$ws.Range("I71").Value = 2589.7
$ws.Range("K71").Value = 12948.5
$ws.Range("M71").Value = -9204.5
$ws.Range("H113").Value = 5155.857
$ws.Range("I113").Value = 7166.3335
$ws.Range("J113").Value = 3648
$ws.Range("K113").Value = 7166.3335
$ws.Range("L113").Value = 3648
$ws.Range("M113").Value = -4996.3335
$ws.Range("N113").Value = -7988
$ws.Range("H122").Value = 7749.5
$ws.Range("I122").Value = 7999.5
$ws.Range("J122").Value = 7499.5
$ws.Range("K122").Value = 23998.5
$ws.Range("L122").Value = 22498.5
$ws.Range("M122").Value = -21548.5
$ws.Range("N122").Value = -27398.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = ""
$ws.Range("H96").Value = 785
$ws.Range("I96").Value = 785
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 785
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 588
$ws.Range("N96").Value = ""
$ws.Range("H107").Value = 877
$ws.Range("I107").Value = 566.7857
$ws.Range("J107").Value = 1497.4286
$ws.Range("K107").Value = 1700.3571
$ws.Range("L107").Value = 4492.2858
$ws.Range("M107").Value = 219.6428999999998
$ws.Range("N107").Value = -8332.2858
